$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old "test21" environment values with the new "test18" ones
$ws.Range("A2").Value = "https://test18.cliotest.com/backoffice/control/main"
$ws.Range("C2").Value = "https://test18.cliotest.com/cabicentral/control/main"
$ws.Range("D2").Value = "https://test18.cliotest.com/warehouse/control/main"
$ws.Range("F2").Value = "virtual_cabitest18"
$ws.Range("G2").Value = "test18"
$ws.Range("K2").Value = "test18"

# Update the active selection in the sheet view
$ws.Range("D14").Select()
